# DEAN import format update
# - Renames "Enrollment" sheet to "Student_Enrollment" and makes it the active tab
# - Adds course_id / delivery_method / campus_name columns to Course_Section
# - Populates the Student_Enrollment header row + a sample data row
# - Updates the selection/active-cell on each sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Course_Section: insert "course_id" and "delivery_method" columns after
# course_section_code, fill in the sample row, and append a "campus_name"
# column at the end.
# ---------------------------------------------------------------------------
$wsCS = $wb.Worksheets.Item("Course_Section")

$wsCS.Range("C1:D1").EntireColumn.Insert()

$wsCS.Cells.Item(1,3).Value = "course_id"
$wsCS.Cells.Item(1,4).Value = "delivery_method"
$wsCS.Cells.Item(2,3).Value = 1395954
$wsCS.Cells.Item(2,4).Value = "online"

# the column insert copies the left neighbour's style (s=3) into C2:D2, but
# the target file leaves those two cells with the default style - borrow it
# from the (still default-styled) instructor_id cell.
$wsCS.Range("I2").Copy()
$wsCS.Range("C2:D2").PasteSpecial(-4122)

# course_section_code sample value changes from "PSY180" to "A01"
$wsCS.Cells.Item(2,2).Value = "A01"

# new trailing campus_name column
$wsCS.Cells.Item(1,10).Value = "campus_name"
$wsCS.Cells.Item(2,10).Value = "Argosy University – Online"

# give the new header cell the same bold header style as the rest of row 1
$wsCS.Range("A1").Copy()
$wsCS.Range("J1").PasteSpecial(-4122)

$wsCS.Columns.Item(4).ColumnWidth = 16.06
$wsCS.Columns.Item(10).ColumnWidth = 22.01
$wsCS.Columns.Item(11).ColumnWidth = 16.06
$wsCS.Rows.Item(1).RowHeight = 15
$wsCS.Rows.Item(2).RowHeight = 15

# ---------------------------------------------------------------------------
# Enrollment -> Student_Enrollment: rename, populate header + sample row,
# and make it the active sheet/tab.
# ---------------------------------------------------------------------------
$wsSE = $wb.Worksheets.Item("Enrollment")
$wsSE.Name = "Student_Enrollment"

# header row formatting matches the other sheets' bold header style (s=2)
$wsCS.Range("A1").Copy()
$wsSE.Range("A1:M1").PasteSpecial(-4122)

$wsSE.Cells.Item(2,1).Value = 49610588
$wsSE.Cells.Item(2,2).Value = 1395954
$wsSE.Cells.Item(2,3).Value = 23849503
$wsSE.Cells.Item(2,6).Value = 1

# course_section_id sample cell reuses the "left" text style (s=3)
$wsCS.Range("A2").Copy()
$wsSE.Range("B2").PasteSpecial(-4122)

$wsSE.Rows.Item(1).RowHeight = 15
$wsSE.Rows.Item(2).RowHeight = 15

$wsSE.Columns.Item(2).ColumnWidth = 16.97
$wsSE.Columns.Item(4).ColumnWidth = 15.97
$wsSE.Columns.Item(5).ColumnWidth = 15.36
$wsSE.Columns.Item(6).ColumnWidth = 23.23
$wsSE.Columns.Item(7).ColumnWidth = 15.46
$wsSE.Columns.Item(8).ColumnWidth = 9.51
$wsSE.Columns.Item(9).ColumnWidth = 29.89
$wsSE.Columns.Item(10).ColumnWidth = 25.74
$wsSE.Columns.Item(11).ColumnWidth = 23.94
$wsSE.Columns.Item(12).ColumnWidth = 28.67
$wsSE.Columns.Item(13).ColumnWidth = 26.86

# ---------------------------------------------------------------------------
# Person / Course: just a selection change. NOTE: Range.Select() activates
# its sheet as a side effect, so these run *before* the final
# Student_Enrollment activation below to make sure that one "wins" and ends
# up as the active tab.
# ---------------------------------------------------------------------------
$wsPerson = $wb.Worksheets.Item("Person")
$wsPerson.Range("A2").Select()

$wsCourse = $wb.Worksheets.Item("Course")
$wsCourse.Range("B2").Select()

$wsCS.Range("B3").Select()

$wsSE.Activate()
$wsSE.Range("F3").Select()

Write-Output "DEAN import format updated"
